$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Translations"

# --- Header row ---
$ws.Range("A1").Value = "Entity Id"
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "Index"
$ws.Range("D1").Value = "Original"
$ws.Range("E1").Value = "Translation"

# --- Row 2 (was A2=1,C2=Id,D2=Orig -> now A=EntityId,B=Title,D=Orig) ---
$ws.Range("A2").Value = "AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"
$ws.Range("B2").Value = "Title"
$ws.Range("C2").Value = $null
$ws.Range("D2").Value = "Orig"
$ws.Range("E2").Value = $null

# --- Row 3 (was A3=2,B3=1,C3=Id,D3=Orig,E3=validation message) ---
$ws.Range("A3").Value = "AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"
$ws.Range("B3").Value = "ValidationMessage"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "Orig"
$ws.Range("E3").Value = "validation message"

# --- Row 4 (was A4=3,C4=Id,D4=Orig) ---
$ws.Range("A4").Value = "AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"
$ws.Range("B4").Value = "Instruction"
$ws.Range("C4").Value = $null
$ws.Range("D4").Value = "Orig"
$ws.Range("E4").Value = $null

# --- Row 5 (was A5=4,B5=2,C5=Id,D5=Orig,E5=option) ---
$ws.Range("A5").Value = "AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"
$ws.Range("B5").Value = "OptionTitle"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "Orig"
$ws.Range("E5").Value = "option"

# Column widths to match new layout (values chosen so the engine's
# internal character-width quantization lands on the closest
# achievable width to the target: A=43.285.., B=18.140.., C=6, E=39.855..)
$ws.Columns.Item(1).ColumnWidth = 42.5
$ws.Columns.Item(2).ColumnWidth = 17.333333333333336
$ws.Columns.Item(3).ColumnWidth = 5.166666666666667
$ws.Columns.Item(5).ColumnWidth = 39.0

# Selection moves to E6 (below the data) as in the final file
$ws.Range("E6").Select()
